$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.797.93'
$ws.Range('E2').Value = '  -1.35%  '
$ws.Range('D3').Value = '2.213.52'
$ws.Range('E3').Value = '  -1.42%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '241.10'
$ws.Range('E5').Value = '  -2.02%  '
$ws.Range('D7').Value = '72.71'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').Value = '0.598'
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').Value = '41.77'
$ws.Range('E10').Value = '  -5.04%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('E12').Value = '  -3.76%  '
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '2.545.86'
$ws.Range('E14').Value = '  -1.39%  '
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').Value = '2.224.15'
$ws.Range('E17').Value = '  -0.64%  '
$ws.Range('D18').Value = '41.696.61'
$ws.Range('E18').Value = '  -1.22%  '
$ws.Range('D19').Value = '0.0000104'
$ws.Range('E19').Value = '  +2.57%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -0.53%  '
$ws.Range('D22').Value = '11.04'
$ws.Range('E22').Value = '  +20.86%  '
$ws.Range('D23').Value = '228.88'
$ws.Range('E23').Value = '  -1.20%  '
$ws.Range('E24').Value = '  -8.72%  '
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('D26').Value = '11.33'
$ws.Range('E26').Value = '  -1.08%  '
$ws.Range('D27').Value = '3.62'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('E28').Value = '  -1.68%  '
$ws.Range('E29').Value = '  -0.72%  '
$ws.Range('D30').Value = '166.96'
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('D31').Value = '20.36'
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('D32').Value = '5.58'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('E33').Value = '  -3.53%  '
$ws.Range('D34').Value = '30.14'
$ws.Range('E34').Value = '  -2.48%  '
$ws.Range('D35').Value = '0.123'
$ws.Range('E35').Value = '  -1.65%  '
$ws.Range('E36').Value = '  -11.01%  '
$ws.Range('E37').Value = '  -6.81%  '
$ws.Range('E38').Value = '  -4.96%  '
$ws.Range('D39').Value = '13.48'
$ws.Range('E39').Value = '  -3.31%  '
$ws.Range('E40').Value = '  -2.89%  '
$ws.Range('D41').Value = '5.61'
$ws.Range('E41').Value = '  -3.22%  '
$ws.Range('D42').Value = '63.78'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('E44').Value = '  -1.63%  '
$ws.Range('D45').Value = '102.75'
$ws.Range('E45').Value = '  -4.61%  '
$ws.Range('D46').Value = '0.0996'
$ws.Range('E46').Value = '  -3.11%  '
$ws.Range('D47').Value = '2.32'
$ws.Range('E47').Value = '  -1.22%  '
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('E49').Value = '  -1.93%  '
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '2.420.22'
$ws.Range('E51').Value = '  -1.56%  '
